# Auto-generated edit script: updates crypto price/volume table
# to match the refreshed data snapshot (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.228.87'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').Value = '3.117.07'
$ws.Range('E3').Value = '  +0.30%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '''579.91'
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('D6').Value = '''173.50'
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '''0.521'
$ws.Range('E8').Value = '  -0.58%  '
$ws.Range('D9').Value = '''6.51'
$ws.Range('E9').Value = '  +1.04%  '
$ws.Range('E10').Value = '  -0.83%  '
$ws.Range('D11').Value = '''0.480'
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('D13').Value = '''36.92'
$ws.Range('E13').Value = '  -1.03%  '
$ws.Range('E14').Value = '  -1.68%  '
$ws.Range('D15').Value = '3.626.54'
$ws.Range('E15').Value = '  +0.15%  '
$ws.Range('D16').Value = '67.177.07'
$ws.Range('E16').Value = '  +0.06%  '
$ws.Range('D17').Value = '''7.10'
$ws.Range('E17').Value = '  -1.46%  '
$ws.Range('D18').Value = '3.111.20'
$ws.Range('E18').Value = '  +0.22%  '
$ws.Range('E19').Value = '  +1.38%  '
$ws.Range('D20').Value = '''490.80'
$ws.Range('E20').Value = '  +1.07%  '
$ws.Range('E21').Value = '  +4.74%  '
$ws.Range('D22').Value = '''0.706'
$ws.Range('E22').Value = '  -1.77%  '
$ws.Range('D23').Value = '''83.99'
$ws.Range('E23').Value = '  -0.29%  '
$ws.Range('D24').Value = '''13.20'
$ws.Range('E24').Value = '  -0.42%  '
$ws.Range('D25').Value = '''2.30'
$ws.Range('E25').Value = '  -3.25%  '
$ws.Range('D26').Value = '''10.54'
$ws.Range('E26').Value = '  +4.91%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').Value = '''7.95'
$ws.Range('E28').Value = '  -0.86%  '
$ws.Range('D29').Value = '''2.38'
$ws.Range('E29').Value = '  -1.14%  '
$ws.Range('D30').Value = '''2.67'
$ws.Range('E30').Value = '  -0.42%  '
$ws.Range('D31').Value = '''28.46'
$ws.Range('E31').Value = '  -1.37%  '
$ws.Range('E32').Value = '  -0.49%  '
$ws.Range('D33').Value = '0.0₃0947'
$ws.Range('E33').Value = '  -6.15%  '
$ws.Range('D34').Value = '''0.998'
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('D35').Value = '''5.88'
$ws.Range('E35').Value = '  -0.61%  '
$ws.Range('D36').Value = '''0.975'
$ws.Range('E36').Value = '  -1.70%  '
$ws.Range('D37').Value = '''47.11'
$ws.Range('E37').Value = '  -1.88%  '
$ws.Range('D38').Value = '''2.05'
$ws.Range('D39').Value = '''0.310'
$ws.Range('E39').Value = '  -2.05%  '
$ws.Range('E40').Value = '  +0.81%  '
$ws.Range('D41').Value = '''8.47'
$ws.Range('E41').Value = '  -2.47%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').Value = '''386.32'
$ws.Range('E42').Value = '  +0.30%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.808.75'
$ws.Range('E43').Value = '  -1.65%  '
$ws.Range('D44').Value = '''2.60'
$ws.Range('E44').Value = '  -7.80%  '
$ws.Range('D45').Value = '''0.0352'
$ws.Range('E45').Value = '  -2.47%  '
$ws.Range('D46').Value = '''135.44'
$ws.Range('E46').Value = '  -0.39%  '
$ws.Range('D48').Value = '''24.96'
$ws.Range('E48').Value = '  -0.33%  '
$ws.Range('D49').Value = '''2.21'
$ws.Range('E49').Value = '  -1.47%  '
$ws.Range('D50').Value = '''0.108'
$ws.Range('E50').Value = '  -1.23%  '
$ws.Range('D51').Value = '''6.73'
$ws.Range('E51').Value = '  -2.19%  '

# Cells that received a leading apostrophe to stop Excel's numeric
# auto-conversion keep the text value but pick up a quotePrefix style;
# clear formatting on just those cells so the style stays untouched.
$ws.Range('D4').ClearFormats()
$ws.Range('D5').ClearFormats()
$ws.Range('D6').ClearFormats()
$ws.Range('D8').ClearFormats()
$ws.Range('D9').ClearFormats()
$ws.Range('D11').ClearFormats()
$ws.Range('D13').ClearFormats()
$ws.Range('D17').ClearFormats()
$ws.Range('D20').ClearFormats()
$ws.Range('D22').ClearFormats()
$ws.Range('D23').ClearFormats()
$ws.Range('D24').ClearFormats()
$ws.Range('D25').ClearFormats()
$ws.Range('D26').ClearFormats()
$ws.Range('D28').ClearFormats()
$ws.Range('D29').ClearFormats()
$ws.Range('D30').ClearFormats()
$ws.Range('D31').ClearFormats()
$ws.Range('D34').ClearFormats()
$ws.Range('D35').ClearFormats()
$ws.Range('D36').ClearFormats()
$ws.Range('D37').ClearFormats()
$ws.Range('D38').ClearFormats()
$ws.Range('D39').ClearFormats()
$ws.Range('D41').ClearFormats()
$ws.Range('D42').ClearFormats()
$ws.Range('D44').ClearFormats()
$ws.Range('D45').ClearFormats()
$ws.Range('D46').ClearFormats()
$ws.Range('D48').ClearFormats()
$ws.Range('D49').ClearFormats()
$ws.Range('D50').ClearFormats()
$ws.Range('D51').ClearFormats()
